# Corrected IFRS financial figures for 엔에스쇼핑 (rows 2-9, columns D:AJ).
# Values were re-keyed from the source feed (prior figures were off by several
# orders of magnitude / pulled the wrong columns); some columns that no longer
# apply to this filing (J, O, and a handful of per-row ratio cells) are cleared.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rowData = @{
    2 = @{ "D" = 3925; "E" = 916; "F" = 916; "G" = 938; "H" = 710; "I" = 710; "J" = $null; "K" = 4153; "L" = 1457; "M" = 2696; "N" = 2696; "O" = $null; "P" = 168; "Q" = 1197; "R" = -650; "S" = -261; "T" = 231; "U" = 966; "V" = 12; "W" = 23.34; "X" = 18.08; "Y" = $null; "Z" = $null; "AA" = 54.04; "AB" = 1490.69; "AC" = 2106; "AD" = $null; "AE" = 8002; "AF" = 0; "AG" = 100; "AH" = $null; "AI" = 4.75; "AJ" = 33696000 }
    3 = @{ "D" = 4064; "E" = 899; "F" = 899; "G" = 928; "H" = 677; "I" = 677; "J" = $null; "K" = 4498; "L" = 1161; "M" = 3338; "N" = 3338; "O" = $null; "P" = 168; "Q" = 435; "R" = 136; "S" = -47; "T" = 225; "U" = 211; "V" = $null; "W" = 22.11; "X" = 16.65; "Y" = 22.43; "Z" = 15.65; "AA" = 34.77; "AB" = 1867.86; "AC" = 2008; "AD" = 9.24; "AE" = 9907; "AF" = 1.87; "AG" = 0; "AH" = 0; "AI" = 0; "AJ" = 33696000 }
    4 = @{ "D" = 4411; "E" = 790; "F" = 790; "G" = 733; "H" = 503; "I" = 503; "J" = $null; "K" = 8087; "L" = 4301; "M" = 3786; "N" = 3786; "O" = $null; "P" = 168; "Q" = 653; "R" = -4858; "S" = 3087; "T" = 4879; "U" = -4226; "V" = 3095; "W" = 17.91; "X" = 11.4; "Y" = 14.12; "Z" = 7.99; "AA" = 113.61; "AB" = 2164.07; "AC" = 1493; "AD" = 10.28; "AE" = 11237; "AF" = 1.37; "AG" = 200; "AH" = 1.3; "AI" = 13.39; "AJ" = 33696000 }
    5 = @{ "D" = 4768; "E" = 800; "F" = 800; "G" = 812; "H" = 552; "I" = 552; "J" = $null; "K" = 7686; "L" = 3517; "M" = 4169; "N" = 4169; "O" = $null; "P" = 168; "Q" = 606; "R" = -15; "S" = -1038; "T" = 157; "U" = 448; "V" = 2226; "W" = 16.78; "X" = 11.57; "Y" = 13.88; "Z" = 7; "AA" = 84.37; "AB" = 2450.24; "AC" = 1638; "AD" = 9.949999999999999; "AE" = 12604; "AF" = 1.29; "AG" = 200; "AH" = 1.23; "AI" = 11.99; "AJ" = 33696000 }
    6 = @{ "D" = 4741; "E" = 612; "F" = 612; "G" = 523; "H" = 346; "I" = 346; "K" = 8140; "L" = 3728; "M" = 4412; "N" = 4412; "P" = 168; "Q" = 298; "R" = -380; "S" = 324; "T" = 306; "U" = -8; "V" = 2623; "W" = 12.9; "X" = 7.29; "Y" = 8.06; "Z" = 4.37; "AA" = 84.48999999999999; "AB" = 2610.92; "AC" = 1026; "AD" = 12.48; "AE" = 13337; "AF" = 0.96; "AG" = 150; "AH" = 1.17; "AI" = 14.36; "AJ" = 33696000 }
    7 = @{ "D" = 4854; "E" = 358; "G" = 149; "H" = 28; "I" = 25; "K" = 8044; "L" = 3672; "M" = 4372; "N" = 4369; "P" = 168; "Q" = 116; "R" = -300; "S" = -45; "T" = 153; "U" = -140; "W" = 7.36; "X" = 0.59; "Y" = 0.57; "Z" = 0.35; "AA" = 83.98; "AC" = 74; "AD" = 118.21; "AE" = 13207; "AF" = 0.66; "AG" = 150; "AH" = 1.71; "AI" = 202.18 }
    8 = @{ "D" = 5296; "E" = 312; "G" = 230; "H" = 167; "I" = 164; "K" = 8092; "L" = 3621; "M" = 4470; "N" = 4464; "P" = 168; "Q" = 236; "R" = -154; "S" = -51; "T" = 125; "U" = 14; "W" = 5.88; "X" = 3.15; "Y" = 3.7; "Z" = 2.07; "AA" = 81; "AC" = 485; "AD" = 18.07; "AE" = 13494; "AF" = 0.65; "AG" = 150; "AH" = 1.71; "AI" = 30.91 }
    9 = @{ "D" = 5786; "E" = 334; "G" = 254; "H" = 182; "I" = 178; "K" = 8160; "L" = 3576; "M" = 4584; "N" = 4574; "P" = 168; "Q" = 206; "R" = -157; "S" = -51; "T" = 125; "U" = 26; "W" = 5.77; "X" = 3.14; "Y" = 3.95; "Z" = 2.23; "AA" = 78.02; "AC" = 530; "AD" = 16.56; "AE" = 13828; "AF" = 0.63; "AG" = 150; "AH" = 1.71; "AI" = 28.32 }
}

foreach ($rowNum in $rowData.Keys) {
    $cellMap = $rowData[$rowNum]
    foreach ($col in $cellMap.Keys) {
        $cellRef = "$col$rowNum"
        $newVal = $cellMap[$col]
        if ($null -eq $newVal) {
            $ws.Range($cellRef).ClearContents()
        } else {
            $ws.Range($cellRef).Value = $newVal
        }
    }
}

